# Insert a new data row at sheet row 65, shifting the existing rows
# 65-140 down to 66-141 (dimension grows from A1:R140 to A1:R141).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value = 8
$ws.Cells.Item(65, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(65, 3).Value = "Coquimbo"
$ws.Cells.Item(65, 4).Value = 44902
$ws.Cells.Item(65, 5).Value = 4
$ws.Cells.Item(65, 6).Value = 100112052
$ws.Cells.Item(65, 7).Value = "Albahaca"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 1060
$ws.Cells.Item(65, 11).Value = 4000
$ws.Cells.Item(65, 12).Value = 4500
$ws.Cells.Item(65, 13).Value = 4250
$ws.Cells.Item(65, 14).Value = "$/paquete"
$ws.Cells.Item(65, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(65, 16).Value = 4250
$ws.Cells.Item(65, 17).Value = 1
$ws.Cells.Item(65, 18).Value = "Hortaliza"
